$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 5425
$ws.Range("I52").Value = 850
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 2550
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -2390
$ws.Range("N52").Value = -30320
$ws.Range("H74").Value = 4776.923
$ws.Range("I74").Value = 4655.5557
$ws.Range("J74").Value = 5050
$ws.Range("K74").Value = 4655.5557
$ws.Range("L74").Value = 5050
$ws.Range("M74").Value = -3719.5557
$ws.Range("N74").Value = -6922
$ws.Range("H76").Value = 4153.3335
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 4153.3335
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 4153.3335
$ws.Range("M76").Value = $null
$ws.Range("N76").Value = -4783.3335
$ws.Range("H77").Value = 4776.923
$ws.Range("I77").Value = 4655.5557
$ws.Range("J77").Value = 5050
$ws.Range("K77").Value = 23277.7785
$ws.Range("L77").Value = 25250
$ws.Range("M77").Value = -18597.7785
$ws.Range("N77").Value = -34610
$ws.Range("H79").Value = 4153.3335
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 4153.3335
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 4153.3335
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -6337.3335
$ws.Range("H128").Value = 250028240
$ws.Range("J128").Value = 250028240
$ws.Range("L128").Value = 250028240
$ws.Range("N128").Value = -250038200
$ws.Range("H130").Value = 333359260
$ws.Range("I130").Value = 38780
$ws.Range("J130").Value = 500019500
$ws.Range("K130").Value = 38780
$ws.Range("L130").Value = 500019500
$ws.Range("M130").Value = -33760
$ws.Range("N130").Value = -500029540
$ws.Range("H137").Value = 1477.5714
$ws.Range("I137").Value = 1923.0869
$ws.Range("J137").Value = 1083.4615
$ws.Range("K137").Value = 5769.2607
$ws.Range("L137").Value = 3250.3845
$ws.Range("M137").Value = -3219.2607
$ws.Range("N137").Value = -8350.3845
$ws.Range("H141").Value = 8245.956
$ws.Range("I141").Value = 3609.9333
$ws.Range("J141").Value = 16938.5
$ws.Range("K141").Value = 10829.7999
$ws.Range("L141").Value = 50815.5
$ws.Range("M141").Value = -5649.7999
$ws.Range("N141").Value = -61175.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6619.5186
$ws.Range("I32").Value = 6306.2466
$ws.Range("J32").Value = 12650
$ws.Range("K32").Value = 6306.2466
$ws.Range("L32").Value = 12650
$ws.Range("M32").Value = -6019.2466
$ws.Range("N32").Value = -13224
$ws.Range("H45").Value = 1537.1111
$ws.Range("I45").Value = 1165.75
$ws.Range("J45").Value = 1834.2
$ws.Range("K45").Value = 1165.75
$ws.Range("L45").Value = 1834.2
$ws.Range("M45").Value = -788.75
$ws.Range("N45").Value = -2588.2
$ws.Range("H61").Value = 5980.5557
$ws.Range("I61").Value = 13956
$ws.Range("J61").Value = 1992.8334
$ws.Range("K61").Value = 13956
$ws.Range("L61").Value = 1992.8334
$ws.Range("M61").Value = -13744
$ws.Range("N61").Value = -2416.8334
$ws.Range("H63").Value = 3803.7896
$ws.Range("I63").Value = 2244.5833
$ws.Range("J63").Value = 6476.7144
$ws.Range("K63").Value = 2244.5833
$ws.Range("L63").Value = 6476.7144
$ws.Range("M63").Value = -1558.5833
$ws.Range("N63").Value = -7848.7144
$ws.Range("H66").Value = 3803.7896
$ws.Range("I66").Value = 2244.5833
$ws.Range("J66").Value = 6476.7144
$ws.Range("K66").Value = 11222.9165
$ws.Range("L66").Value = 32383.572
$ws.Range("M66").Value = -7790.916499999999
$ws.Range("N66").Value = -39247.572
$ws.Range("H132").Value = 870132.9399999999
$ws.Range("I132").Value = 2270571.2
$ws.Range("J132").Value = 5156.4116
$ws.Range("K132").Value = 6811713.600000001
$ws.Range("L132").Value = 15469.2348
$ws.Range("M132").Value = -6809183.600000001
$ws.Range("N132").Value = -20529.2348
$ws.Range("H133").Value = 37666.668
$ws.Range("J133").Value = 37666.668
$ws.Range("L133").Value = 37666.668
$ws.Range("N133").Value = -42726.668
$ws.Range("H136").Value = 5980.5557
$ws.Range("I136").Value = 13956
$ws.Range("J136").Value = 1992.8334
$ws.Range("K136").Value = 41868
$ws.Range("L136").Value = 5978.5002
$ws.Range("M136").Value = -39318
$ws.Range("N136").Value = -11078.5002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1555.6923
$ws.Range("I22").Value = 1874.8889
$ws.Range("J22").Value = 837.5
$ws.Range("K22").Value = 1874.8889
$ws.Range("L22").Value = 837.5
$ws.Range("M22").Value = -1701.8889
$ws.Range("N22").Value = -1183.5
$ws.Range("H105").Value = 1819.875
$ws.Range("J105").Value = 2925
$ws.Range("L105").Value = 2925
$ws.Range("N105").Value = -6419
$ws.Range("H134").Value = 5921.974
$ws.Range("I134").Value = 2668.0833
$ws.Range("J134").Value = 7368.148
$ws.Range("K134").Value = 8004.249899999999
$ws.Range("L134").Value = 22104.444
$ws.Range("M134").Value = -5469.249899999999
$ws.Range("N134").Value = -27174.444
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7564.1577
$ws.Range("I62").Value = 2222.7856
$ws.Range("J62").Value = 22520
$ws.Range("K62").Value = 2222.7856
$ws.Range("L62").Value = 22520
$ws.Range("M62").Value = -1598.7856
$ws.Range("N62").Value = -23768
$ws.Range("H65").Value = 7564.1577
$ws.Range("I65").Value = 2222.7856
$ws.Range("J65").Value = 22520
$ws.Range("K65").Value = 11113.928
$ws.Range("L65").Value = 112600
$ws.Range("M65").Value = -7993.928
$ws.Range("N65").Value = -118840
$ws.Range("H132").Value = 3400.5557
$ws.Range("I132").Value = 3589.2222
$ws.Range("K132").Value = 10767.6666
$ws.Range("M132").Value = -8237.6666
$ws.Range("H134").Value = 2300.077
$ws.Range("I134").Value = 769
$ws.Range("J134").Value = 2759.4
$ws.Range("K134").Value = 2307
$ws.Range("L134").Value = 8278.200000000001
$ws.Range("M134").Value = 228
$ws.Range("N134").Value = -13348.2
$ws.Range("H135").Value = 33540
$ws.Range("J135").Value = 33540
$ws.Range("L135").Value = 33540
$ws.Range("N135").Value = -43680
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 838.3333
$ws.Range("I51").Value = 838.3333
$ws.Range("K51").Value = 2514.9999
$ws.Range("M51").Value = -2054.9999
$ws.Range("H110").Value = 2304.5715
$ws.Range("I110").Value = 2304.5715
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 6913.7145
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -2823.7145
$ws.Range("N110").Value = $null
$ws.Range("H120").Value = 19999
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").Value = $null
$ws.Range("H131").Value = 767.84
$ws.Range("I131").Value = 459.86667
$ws.Range("J131").Value = 1229.8
$ws.Range("K131").Value = 1379.60001
$ws.Range("L131").Value = 3689.4
$ws.Range("M131").Value = 3660.39999
$ws.Range("N131").Value = -13769.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 52203.6
$ws.Range("J20").Value = 65004.5
$ws.Range("L20").Value = 65004.5
$ws.Range("N20").Value = -65494.5
$ws.Range("H70").Value = 6318
$ws.Range("I70").Value = 5707.8
$ws.Range("J70").Value = 6787.385
$ws.Range("K70").Value = 5707.8
$ws.Range("L70").Value = 6787.385
$ws.Range("M70").Value = -5437.8
$ws.Range("N70").Value = -7327.385
$ws.Range("H73").Value = 6318
$ws.Range("I73").Value = 5707.8
$ws.Range("J73").Value = 6787.385
$ws.Range("K73").Value = 5707.8
$ws.Range("L73").Value = 6787.385
$ws.Range("M73").Value = -4771.8
$ws.Range("N73").Value = -8659.385
$ws.Range("H80").Value = 3047.5
$ws.Range("I80").Value = 2657
$ws.Range("K80").Value = 2657
$ws.Range("M80").Value = -1659
$ws.Range("H83").Value = 3047.5
$ws.Range("I83").Value = 2657
$ws.Range("K83").Value = 13285
$ws.Range("M83").Value = -8293
$ws.Range("H122").Value = 5093.4546
$ws.Range("I122").Value = 4503.5
$ws.Range("J122").Value = 6666.6665
$ws.Range("K122").Value = 13510.5
$ws.Range("L122").Value = 19999.9995
$ws.Range("M122").Value = -11060.5
$ws.Range("N122").Value = -24899.9995
$ws.Range("H132").Value = 2196220
$ws.Range("I132").Value = 10419916
$ws.Range("J132").Value = 3234.1333
$ws.Range("K132").Value = 31259748
$ws.Range("L132").Value = 9702.3999
$ws.Range("M132").Value = -31257218
$ws.Range("N132").Value = -14762.3999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 79943
$ws.Range("J46").Value = 79943
$ws.Range("L46").Value = 79943
$ws.Range("N46").Value = -80405
$ws.Range("H134").Value = 79943
$ws.Range("J134").Value = 79943
$ws.Range("L134").Value = 239829
$ws.Range("N134").Value = -244899
